$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3070.5715
$ws.Range("I19").Value = 3000
$ws.Range("J19").Value = 3082.3333
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3082.3333
$ws.Range("M19").Value = -2825
$ws.Range("N19").Value = -3432.3333

$ws.Range("H98").Value = 111113330
$ws.Range("I98").Value = 111113330
$ws.Range("K98").Value = 111113330
$ws.Range("M98").Value = -111111832

$ws.Range("H107").Value = 404.6316
$ws.Range("I107").Value = 421.6111
$ws.Range("K107").Value = 421.6111
$ws.Range("M107").Value = 1498.3889

$ws.Range("H113").Value = 100002150
$ws.Range("I113").Value = 50000376
$ws.Range("J113").Value = 133336664
$ws.Range("K113").Value = 50000376
$ws.Range("L113").Value = 133336664
$ws.Range("M113").Value = -49997122
$ws.Range("N113").Value = -133343172

$ws.Range("H122").Value = 111113330
$ws.Range("I122").Value = 111113330
$ws.Range("K122").Value = 333339990
$ws.Range("M122").Value = -333337540

$ws.Range("H131").Value = 3247.6875
$ws.Range("I131").Value = 1087.5454
$ws.Range("K131").Value = 3262.6362
$ws.Range("M131").Value = 1777.3638

$ws.Range("H135").Value = 2055.6191
$ws.Range("I135").Value = 2252.5625
$ws.Range("J135").Value = 1425.4
$ws.Range("K135").Value = 20273.0625
$ws.Range("L135").Value = 12828.6
$ws.Range("M135").Value = -17738.0625
$ws.Range("N135").Value = -17898.6

$ws.Range("H137").Value = 3365.56
$ws.Range("I137").Value = 3611.9048
$ws.Range("K137").Value = 10835.7144
$ws.Range("M137").Value = -8285.714399999999

$ws.Range("H138").Value = 2251.457
$ws.Range("J138").Value = 2690.86
$ws.Range("L138").Value = 8072.58
$ws.Range("N138").Value = -18352.58

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 3385
$ws.Range("I14").Value = 5002.5
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 5002.5
$ws.Range("L14").Value = 150
$ws.Range("M14").Value = -4827.5
$ws.Range("N14").Value = -500

$ws.Range("H32").Value = 26334966
$ws.Range("I32").Value = 33349512
$ws.Range("J32").Value = 30411.75
$ws.Range("K32").Value = 33349512
$ws.Range("L32").Value = 30411.75
$ws.Range("M32").Value = -33349225
$ws.Range("N32").Value = -30985.75

$ws.Range("H104").Value = 37784.8
$ws.Range("J104").Value = 37784.8
$ws.Range("L104").Value = 37784.8
$ws.Range("N104").Value = -44772.8

$ws.Range("H135").Value = 500050000
$ws.Range("J135").Value = 500050000
$ws.Range("L135").Value = 500050000
$ws.Range("N135").Value = -500060140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4221.7144
$ws.Range("I20").Value = 4511.1
$ws.Range("J20").Value = 3498.25
$ws.Range("K20").Value = 4511.1
$ws.Range("L20").Value = 3498.25
$ws.Range("M20").Value = -4264.1
$ws.Range("N20").Value = -3992.25

$ws.Range("H135").Value = 60238.438
$ws.Range("J135").Value = 60238.438
$ws.Range("L135").Value = 60238.438
$ws.Range("N135").Value = -70378.43799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 548653.75
$ws.Range("J31").Value = 1230590
$ws.Range("L31").Value = 1230590
$ws.Range("N31").Value = -1231180

$ws.Range("H34").Value = 548653.75
$ws.Range("J34").Value = 1230590
$ws.Range("L34").Value = 1230590
$ws.Range("N34").Value = -1230994

$ws.Range("H99").Value = 3999.6
$ws.Range("J99").Value = 4000
$ws.Range("L99").Value = 4000
$ws.Range("N99").Value = -6996

$ws.Range("H126").Value = 3999.6
$ws.Range("J126").Value = 4000
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 814.5
$ws.Range("J34").Value = 1500
$ws.Range("L34").Value = 4500
$ws.Range("N34").Value = -4668

$ws.Range("H114").Value = 665
$ws.Range("I114").Value = 864
$ws.Range("J114").Value = 399.66666
$ws.Range("K114").Value = 2592
$ws.Range("L114").Value = 1198.99998
$ws.Range("M114").Value = 662
$ws.Range("N114").Value = -7706.999980000001

$ws.Range("H122").Value = 2565.9524
$ws.Range("J122").Value = 3663.2144
$ws.Range("L122").Value = 32968.9296
$ws.Range("N122").Value = -37868.9296

$ws.Range("H131").Value = 6157.365
$ws.Range("J131").Value = 5795.5
$ws.Range("L131").Value = 17386.5
$ws.Range("N131").Value = -27466.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 87233.89
$ws.Range("I19").Value = 128350.836
$ws.Range("K19").Value = 128350.836
$ws.Range("M19").Value = -128062.836

$ws.Range("H62").Value = 115000
$ws.Range("J62").Value = 115000
$ws.Range("L62").Value = 115000
$ws.Range("N62").Value = -116372

$ws.Range("H65").Value = 115000
$ws.Range("J65").Value = 115000
$ws.Range("L65").Value = 345000
$ws.Range("N65").Value = -351864

$ws.Range("H110").Value = 96369
$ws.Range("J110").Value = 96369
$ws.Range("L110").Value = 96369
$ws.Range("N110").Value = -104549

$ws.Range("H111").Value = 85000
$ws.Range("J111").Value = 85000
$ws.Range("L111").Value = 85000
$ws.Range("N111").Value = -91134

$ws.Range("H113").Value = 3503.0588
$ws.Range("I113").Value = 2171.3333
$ws.Range("K113").Value = 2171.3333
$ws.Range("M113").Value = -1.333299999999781

$ws.Range("H132").Value = 100013310
$ws.Range("I132").Value = 200002540
$ws.Range("J132").Value = 24081.2
$ws.Range("K132").Value = 600007620
$ws.Range("L132").Value = 72243.60000000001
$ws.Range("M132").Value = -600005090
$ws.Range("N132").Value = -77303.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 125000000
$ws.Range("J44").Value = 125000000
$ws.Range("L44").Value = 125000000
$ws.Range("N44").Value = -125000912

$ws.Range("I55").Value = 30303504
$ws.Range("J55").Value = 621.8889
$ws.Range("K55").Value = 30303504
$ws.Range("L55").Value = 621.8889
$ws.Range("M55").Value = -30303331
$ws.Range("N55").Value = -967.8889

$ws.Range("H58").Value = 7357.143
$ws.Range("I58").Value = 4750
$ws.Range("J58").Value = 8400
$ws.Range("K58").Value = 4750
$ws.Range("L58").Value = 8400
$ws.Range("M58").Value = -4490
$ws.Range("N58").Value = -8920

$ws.Range("H74").Value = 82500
$ws.Range("I74").Value = 82500
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 82500
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -81502
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 82500
$ws.Range("I77").Value = 82500
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 247500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -242508
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1142.8572
$ws.Range("I14").Value = 200
$ws.Range("J14").Value = 3500
$ws.Range("K14").Value = 200
$ws.Range("L14").Value = 3500
$ws.Range("M14").Value = -32
$ws.Range("N14").Value = -3836

$ws.Range("H15").Value = 45003.5
$ws.Range("I15").Value = 10000
$ws.Range("J15").Value = 80007
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 80007
$ws.Range("M15").Value = -9712
$ws.Range("N15").Value = -80583

$ws.Range("H39").Value = 13949.4
$ws.Range("I39").Value = 9875
$ws.Range("J39").Value = 30247
$ws.Range("K39").Value = 9875
$ws.Range("L39").Value = 30247
$ws.Range("M39").Value = -9462
$ws.Range("N39").Value = -31073

$ws.Range("H69").Value = 21666.666
$ws.Range("J69").Value = 21666.666
$ws.Range("L69").Value = 21666.666
$ws.Range("N69").Value = -23164.666

$ws.Range("H72").Value = 21666.666
$ws.Range("J72").Value = 21666.666
$ws.Range("L72").Value = 64999.99800000001
$ws.Range("N72").Value = -72487.99800000001
